# Natmi following Dr Hou advice
#
# The ligand/receptor summary (Fgf13 -> Scn5a) is recomputed after adding
# the "M2" cluster into the sending/target cluster set alongside the
# existing "ECs" and "sCs" clusters. Every cluster-pair row is refreshed
# with the new detection/expression/specificity statistics, and three new
# rows are added for the pairs that involve the "M2" cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression
# value, Ligand total expression value, Ligand derived specificity of
# average expression value, Ligand derived specificity of total expression
# value, Receptor-expressing cells, Receptor detection rate, Receptor
# average expression value, Receptor total expression value, Receptor
# derived specificity of average expression value, Receptor derived
# specificity of total expression value, Edge average expression weight,
# Edge total expression weight, Edge average expression derived
# specificity, Edge total expression derived specificity.
$data = @(
    @("ECs", "Fgf13", "Scn5a", "ECs", 3, 1, 0.1139166666666667, 0.34175, 0.04552287366220362, 0.04552287366220363, 3, 1, 0.48641, 1.45923, 0.1393325029444749, 0.1393325029444749, 0.05541020583333333, 0.4986918524999999, 0.006342815928579945, 0.006342815928579945),
    @("ECs", "Fgf13", "Scn5a", "sCs", 3, 1, 0.1139166666666667, 0.34175, 0.04552287366220362, 0.04552287366220363, 3, 1, 3.004591666666666, 9.013774999999999, 0.8606674970555251, 0.8606674970555251, 0.3422730673611111, 3.08045760625, 0.03918005773362368, 0.03918005773362369),
    @("M2", "Fgf13", "Scn5a", "ECs", 2, 0.6666666666666666, 0.3546683333333333, 1.064005, 0.1417309881227592, 0.1417309881227592, 3, 1, 0.48641, 1.45923, 0.1393325029444749, 0.1393325029444749, 0.1725142240166666, 1.55262801615, 0.01974773331993768, 0.01974773331993768),
    @("M2", "Fgf13", "Scn5a", "sCs", 2, 0.6666666666666666, 0.3546683333333333, 1.064005, 0.1417309881227592, 0.1417309881227592, 3, 1, 3.004591666666666, 9.013774999999999, 0.8606674970555251, 0.8606674970555251, 1.065633518763889, 9.590701668874997, 0.1219832548028215, 0.1219832548028215),
    @("sCs", "Fgf13", "Scn5a", "ECs", 3, 1, 2.03382, 6.101459999999999, 0.8127461382150372, 0.8127461382150372, 3, 1, 0.48641, 1.45923, 0.1393325029444749, 0.1393325029444749, 0.9892703861999999, 8.903433475799998, 0.1132419536959573, 0.1132419536959572),
    @("sCs", "Fgf13", "Scn5a", "sCs", 3, 1, 2.03382, 6.101459999999999, 0.8127461382150372, 0.8127461382150372, 3, 1, 3.004591666666666, 9.013774999999999, 0.8606674970555251, 0.8606674970555251, 6.110798623499999, 54.99718761149999, 0.6995041845190799, 0.6995041845190799)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
